$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "ChildResolutionInput / Letter Sent / Resolution input when resolving the
# child case" row (row 32) is being removed entirely - the parent resolution
# row takes over as the single "ResolutionInput" row.
$ws.Rows.Item(32).Delete()

# Rename what is now row 32 (previously "ParentResolutionInput") to the new
# consolidated name "ResolutionInput".
$ws.Range("A32").Value = "ResolutionInput"

# Resize the worksheet table (Table1) down to the new extent now that a row
# has been removed.
$ws.ListObjects.Item("Table1").Resize($ws.Range("A1:C35"))

# Restore the view: scroll back to the top and select B31 (matches the
# author's final cursor position after the edit).
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("B31").Select()
